$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full corrected data set (rows 2..256) for columns A (name), B (categories), C (amount).
# This replaces the previous data, which had mis-assigned "amount" values for
# several biogenic/non-fossil and other gas entries (an off-by-one style
# mix-up between rows), per "Correct values for biogenic methane."
$data = @(
  @("Bromopropane","air::unspecified",0.188),
  @("Butane","air::urban air close to ground",0.022),
  @("Butane","air::non-urban air or from high stacks",0.022),
  @("Butane","air::low population density, long-term",0.022),
  @("Butane","air::lower stratosphere + upper troposphere",0.022),
  @("Butane","air::unspecified",0.022),
  @("Carbon dioxide, fossil","air::urban air close to ground",1.0),
  @("Carbon dioxide, fossil","air::non-urban air or from high stacks",1.0),
  @("Carbon dioxide, fossil","air::low population density, long-term",1.0),
  @("Carbon dioxide, fossil","air::lower stratosphere + upper troposphere",1.0),
  @("Carbon dioxide, fossil","air::unspecified",1.0),
  @("Carbon dioxide, non-fossil","air::urban air close to ground",1.0),
  @("Carbon dioxide, non-fossil","air::non-urban air or from high stacks",1.0),
  @("Carbon dioxide, non-fossil","air::low population density, long-term",1.0),
  @("Carbon dioxide, non-fossil","air::lower stratosphere + upper troposphere",1.0),
  @("Carbon dioxide, non-fossil","air::unspecified",1.0),
  @("Carbon dioxide, from soil or biomass stock","air::unspecified",1.0),
  @("Carbon dioxide, from soil or biomass stock","air::low population density, long-term",1.0),
  @("Carbon dioxide, from soil or biomass stock","air::lower stratosphere + upper troposphere",1.0),
  @("Carbon dioxide, from soil or biomass stock","air::non-urban air or from high stacks",1.0),
  @("Carbon dioxide, from soil or biomass stock","air::urban air close to ground",1.0),
  @("Carbon dioxide, in air","natural resource::in air",-1.0),
  @("Carbon dioxide, non-fossil, resource correction","natural resource::in air",-1.0),
  @("Carbon dioxide, to soil or biomass stock","soil::unspecified",-1.0),
  @("Carbon dioxide, to soil or biomass stock","soil::agricultural",-1.0),
  @("Carbon dioxide, to soil or biomass stock","soil::industrial",-1.0),
  @("Carbon dioxide, to soil or biomass stock","soil::forestry",-1.0),
  @("Carbon monoxide, fossil","air::urban air close to ground",9.221400000000001),
  @("Carbon monoxide, fossil","air::non-urban air or from high stacks",9.221400000000001),
  @("Carbon monoxide, fossil","air::low population density, long-term",9.221400000000001),
  @("Carbon monoxide, fossil","air::lower stratosphere + upper troposphere",9.221400000000001),
  @("Carbon monoxide, fossil","air::unspecified",9.221400000000001),
  @("Carbon monoxide, from soil or biomass stock","air::low population density, long-term",9.221400000000001),
  @("Carbon monoxide, from soil or biomass stock","air::unspecified",9.221400000000001),
  @("Carbon monoxide, from soil or biomass stock","air::urban air close to ground",9.221400000000001),
  @("Carbon monoxide, from soil or biomass stock","air::non-urban air or from high stacks",9.221400000000001),
  @("Carbon monoxide, from soil or biomass stock","air::lower stratosphere + upper troposphere",9.221400000000001),
  @("Carbon monoxide, non-fossil","air::urban air close to ground",7.65),
  @("Carbon monoxide, non-fossil","air::non-urban air or from high stacks",7.65),
  @("Carbon monoxide, non-fossil","air::low population density, long-term",7.65),
  @("Carbon monoxide, non-fossil","air::lower stratosphere + upper troposphere",7.65),
  @("Carbon monoxide, non-fossil","air::unspecified",7.65),
  @("Chloroform","air::urban air close to ground",74.2),
  @("Chloroform","air::non-urban air or from high stacks",74.2),
  @("Chloroform","air::low population density, long-term",74.2),
  @("Chloroform","air::lower stratosphere + upper troposphere",74.2),
  @("Chloroform","air::unspecified",74.2),
  @("Dinitrogen monoxide","air::urban air close to ground",273.0),
  @("Dinitrogen monoxide","air::non-urban air or from high stacks",273.0),
  @("Dinitrogen monoxide","air::low population density, long-term",273.0),
  @("Dinitrogen monoxide","air::lower stratosphere + upper troposphere",273.0),
  @("Dinitrogen monoxide","air::unspecified",273.0),
  @("Ethane","air::urban air close to ground",1.57),
  @("Ethane","air::non-urban air or from high stacks",1.57),
  @("Ethane","air::low population density, long-term",1.57),
  @("Ethane","air::lower stratosphere + upper troposphere",1.57),
  @("Ethane","air::unspecified",1.57),
  @("Ethane, 1-chloro-1,1-difluoro-, HCFC-142b","air::urban air close to ground",5510.0),
  @("Ethane, 1-chloro-1,1-difluoro-, HCFC-142b","air::non-urban air or from high stacks",5510.0),
  @("Ethane, 1-chloro-1,1-difluoro-, HCFC-142b","air::low population density, long-term",5510.0),
  @("Ethane, 1-chloro-1,1-difluoro-, HCFC-142b","air::lower stratosphere + upper troposphere",5510.0),
  @("Ethane, 1-chloro-1,1-difluoro-, HCFC-142b","air::unspecified",5510.0),
  @("Ethane, 1,1-dichloro-1-fluoro-, HCFC-141b","air::urban air close to ground",2710.0),
  @("Ethane, 1,1-dichloro-1-fluoro-, HCFC-141b","air::non-urban air or from high stacks",2710.0),
  @("Ethane, 1,1-dichloro-1-fluoro-, HCFC-141b","air::low population density, long-term",2710.0),
  @("Ethane, 1,1-dichloro-1-fluoro-, HCFC-141b","air::lower stratosphere + upper troposphere",2710.0),
  @("Ethane, 1,1-dichloro-1-fluoro-, HCFC-141b","air::unspecified",2710.0),
  @("Ethane, 1,1-difluoro-, HFC-152a","air::urban air close to ground",591.0),
  @("Ethane, 1,1-difluoro-, HFC-152a","air::non-urban air or from high stacks",591.0),
  @("Ethane, 1,1-difluoro-, HFC-152a","air::low population density, long-term",591.0),
  @("Ethane, 1,1-difluoro-, HFC-152a","air::lower stratosphere + upper troposphere",591.0),
  @("Ethane, 1,1-difluoro-, HFC-152a","air::unspecified",591.0),
  @("Ethane, 1,1,1-trichloro-, HCFC-140","air::urban air close to ground",567.0),
  @("Ethane, 1,1,1-trichloro-, HCFC-140","air::non-urban air or from high stacks",567.0),
  @("Ethane, 1,1,1-trichloro-, HCFC-140","air::low population density, long-term",567.0),
  @("Ethane, 1,1,1-trichloro-, HCFC-140","air::lower stratosphere + upper troposphere",567.0),
  @("Ethane, 1,1,1-trichloro-, HCFC-140","air::unspecified",567.0),
  @("Ethane, 1,1,1-trifluoro-, HFC-143a","air::urban air close to ground",7840.0),
  @("Ethane, 1,1,1-trifluoro-, HFC-143a","air::non-urban air or from high stacks",7840.0),
  @("Ethane, 1,1,1-trifluoro-, HFC-143a","air::low population density, long-term",7840.0),
  @("Ethane, 1,1,1-trifluoro-, HFC-143a","air::lower stratosphere + upper troposphere",7840.0),
  @("Ethane, 1,1,1-trifluoro-, HFC-143a","air::unspecified",7840.0),
  @("Ethane, 1,1,1,2-tetrafluoro-, HFC-134a","air::urban air close to ground",4140.0),
  @("Ethane, 1,1,1,2-tetrafluoro-, HFC-134a","air::non-urban air or from high stacks",4140.0),
  @("Ethane, 1,1,1,2-tetrafluoro-, HFC-134a","air::low population density, long-term",4140.0),
  @("Ethane, 1,1,1,2-tetrafluoro-, HFC-134a","air::lower stratosphere + upper troposphere",4140.0),
  @("Ethane, 1,1,1,2-tetrafluoro-, HFC-134a","air::unspecified",4140.0),
  @("Ethane, 1,1,2-trichloro-1,2,2-trifluoro-, CFC-113","air::low population density, long-term",6860.0),
  @("Ethane, 1,1,2-trichloro-1,2,2-trifluoro-, CFC-113","air::unspecified",6860.0),
  @("Ethane, 1,1,2-trichloro-1,2,2-trifluoro-, CFC-113","air::lower stratosphere + upper troposphere",6860.0),
  @("Ethane, 1,1,2-trichloro-1,2,2-trifluoro-, CFC-113","air::non-urban air or from high stacks",6860.0),
  @("Ethane, 1,1,2-trichloro-1,2,2-trifluoro-, CFC-113","air::urban air close to ground",6860.0),
  @("Ethane, 1,2-dichloro-","air::urban air close to ground",4.68),
  @("Ethane, 1,2-dichloro-","air::non-urban air or from high stacks",4.68),
  @("Ethane, 1,2-dichloro-","air::low population density, long-term",4.68),
  @("Ethane, 1,2-dichloro-","air::lower stratosphere + upper troposphere",4.68),
  @("Ethane, 1,2-dichloro-","air::unspecified",4.68),
  @("Ethane, 1,2-dichloro-1,1,2,2-tetrafluoro-, CFC-114","air::unspecified",8260.0),
  @("Ethane, 1,2-dichloro-1,1,2,2-tetrafluoro-, CFC-114","air::lower stratosphere + upper troposphere",8260.0),
  @("Ethane, 1,2-dichloro-1,1,2,2-tetrafluoro-, CFC-114","air::low population density, long-term",8260.0),
  @("Ethane, 1,2-dichloro-1,1,2,2-tetrafluoro-, CFC-114","air::non-urban air or from high stacks",8260.0),
  @("Ethane, 1,2-dichloro-1,1,2,2-tetrafluoro-, CFC-114","air::urban air close to ground",8260.0),
  @("Ethane, 2-chloro-1,1,1,2-tetrafluoro-, HCFC-124","air::urban air close to ground",2070.0),
  @("Ethane, 2-chloro-1,1,1,2-tetrafluoro-, HCFC-124","air::non-urban air or from high stacks",2070.0),
  @("Ethane, 2-chloro-1,1,1,2-tetrafluoro-, HCFC-124","air::low population density, long-term",2070.0),
  @("Ethane, 2-chloro-1,1,1,2-tetrafluoro-, HCFC-124","air::lower stratosphere + upper troposphere",2070.0),
  @("Ethane, 2-chloro-1,1,1,2-tetrafluoro-, HCFC-124","air::unspecified",2070.0),
  @("Ethane, 2,2-dichloro-1,1,1-trifluoro-, HCFC-123","air::urban air close to ground",325.0),
  @("Ethane, 2,2-dichloro-1,1,1-trifluoro-, HCFC-123","air::non-urban air or from high stacks",325.0),
  @("Ethane, 2,2-dichloro-1,1,1-trifluoro-, HCFC-123","air::low population density, long-term",325.0),
  @("Ethane, 2,2-dichloro-1,1,1-trifluoro-, HCFC-123","air::lower stratosphere + upper troposphere",325.0),
  @("Ethane, 2,2-dichloro-1,1,1-trifluoro-, HCFC-123","air::unspecified",325.0),
  @("Ethane, chloropentafluoro-, CFC-115","air::urban air close to ground",7410.0),
  @("Ethane, chloropentafluoro-, CFC-115","air::non-urban air or from high stacks",7410.0),
  @("Ethane, chloropentafluoro-, CFC-115","air::low population density, long-term",7410.0),
  @("Ethane, chloropentafluoro-, CFC-115","air::lower stratosphere + upper troposphere",7410.0),
  @("Ethane, chloropentafluoro-, CFC-115","air::unspecified",7410.0),
  @("Ethane, hexafluoro-, HFC-116","air::urban air close to ground",8940.0),
  @("Ethane, hexafluoro-, HFC-116","air::non-urban air or from high stacks",8940.0),
  @("Ethane, hexafluoro-, HFC-116","air::low population density, long-term",8940.0),
  @("Ethane, hexafluoro-, HFC-116","air::lower stratosphere + upper troposphere",8940.0),
  @("Ethane, hexafluoro-, HFC-116","air::unspecified",8940.0),
  @("Ethane, pentafluoro-, HFC-125","air::urban air close to ground",6740.0),
  @("Ethane, pentafluoro-, HFC-125","air::non-urban air or from high stacks",6740.0),
  @("Ethane, pentafluoro-, HFC-125","air::low population density, long-term",6740.0),
  @("Ethane, pentafluoro-, HFC-125","air::lower stratosphere + upper troposphere",6740.0),
  @("Ethane, pentafluoro-, HFC-125","air::unspecified",6740.0),
  @("Hydrogen","air::low population density, long-term",33.0),
  @("Hydrogen","air::lower stratosphere + upper troposphere",33.0),
  @("Hydrogen","air::non-urban air or from high stacks",33.0),
  @("Hydrogen","air::unspecified",33.0),
  @("Hydrogen","air::urban air close to ground",33.0),
  @("Methane, bromo-, Halon 1001","air::urban air close to ground",8.74),
  @("Methane, bromo-, Halon 1001","air::non-urban air or from high stacks",8.74),
  @("Methane, bromo-, Halon 1001","air::low population density, long-term",8.74),
  @("Methane, bromo-, Halon 1001","air::lower stratosphere + upper troposphere",8.74),
  @("Methane, bromo-, Halon 1001","air::unspecified",8.74),
  @("Methane, bromochlorodifluoro-, Halon 1211","air::urban air close to ground",4920.0),
  @("Methane, bromochlorodifluoro-, Halon 1211","air::non-urban air or from high stacks",4920.0),
  @("Methane, bromochlorodifluoro-, Halon 1211","air::low population density, long-term",4920.0),
  @("Methane, bromochlorodifluoro-, Halon 1211","air::lower stratosphere + upper troposphere",4920.0),
  @("Methane, bromochlorodifluoro-, Halon 1211","air::unspecified",4920.0),
  @("Methane, bromotrifluoro-, Halon 1301","air::urban air close to ground",8320.0),
  @("Methane, bromotrifluoro-, Halon 1301","air::non-urban air or from high stacks",8320.0),
  @("Methane, bromotrifluoro-, Halon 1301","air::low population density, long-term",8320.0),
  @("Methane, bromotrifluoro-, Halon 1301","air::lower stratosphere + upper troposphere",8320.0),
  @("Methane, bromotrifluoro-, Halon 1301","air::unspecified",8320.0),
  @("Methane, chloro-fluoro-, HCFC-31","air::urban air close to ground",286.0),
  @("Methane, chloro-fluoro-, HCFC-31","air::non-urban air or from high stacks",286.0),
  @("Methane, chloro-fluoro-, HCFC-31","air::low population density, long-term",286.0),
  @("Methane, chloro-fluoro-, HCFC-31","air::lower stratosphere + upper troposphere",286.0),
  @("Methane, chloro-fluoro-, HCFC-31","air::unspecified",286.0),
  @("Methane, chlorodifluoro-, HCFC-22","air::urban air close to ground",5690.0),
  @("Methane, chlorodifluoro-, HCFC-22","air::non-urban air or from high stacks",5690.0),
  @("Methane, chlorodifluoro-, HCFC-22","air::low population density, long-term",5690.0),
  @("Methane, chlorodifluoro-, HCFC-22","air::lower stratosphere + upper troposphere",5690.0),
  @("Methane, chlorodifluoro-, HCFC-22","air::unspecified",5690.0),
  @("Methane, chlorotrifluoro-, CFC-13","air::urban air close to ground",12400.0),
  @("Methane, chlorotrifluoro-, CFC-13","air::non-urban air or from high stacks",12400.0),
  @("Methane, chlorotrifluoro-, CFC-13","air::low population density, long-term",12400.0),
  @("Methane, chlorotrifluoro-, CFC-13","air::lower stratosphere + upper troposphere",12400.0),
  @("Methane, chlorotrifluoro-, CFC-13","air::unspecified",12400.0),
  @("Methane, dichloro-, HCC-30","air::urban air close to ground",40.2),
  @("Methane, dichloro-, HCC-30","air::non-urban air or from high stacks",40.2),
  @("Methane, dichloro-, HCC-30","air::low population density, long-term",40.2),
  @("Methane, dichloro-, HCC-30","air::lower stratosphere + upper troposphere",40.2),
  @("Methane, dichloro-, HCC-30","air::unspecified",40.2),
  @("Methane, dichlorodifluoro-, CFC-12","air::urban air close to ground",12700.0),
  @("Methane, dichlorodifluoro-, CFC-12","air::non-urban air or from high stacks",12700.0),
  @("Methane, dichlorodifluoro-, CFC-12","air::low population density, long-term",12700.0),
  @("Methane, dichlorodifluoro-, CFC-12","air::lower stratosphere + upper troposphere",12700.0),
  @("Methane, dichlorodifluoro-, CFC-12","air::unspecified",12700.0),
  @("Methane, dichlorofluoro-, HCFC-21","air::urban air close to ground",575.0),
  @("Methane, dichlorofluoro-, HCFC-21","air::non-urban air or from high stacks",575.0),
  @("Methane, dichlorofluoro-, HCFC-21","air::low population density, long-term",575.0),
  @("Methane, dichlorofluoro-, HCFC-21","air::lower stratosphere + upper troposphere",575.0),
  @("Methane, dichlorofluoro-, HCFC-21","air::unspecified",575.0),
  @("Methane, difluoro-, HFC-32","air::urban air close to ground",2690.0),
  @("Methane, difluoro-, HFC-32","air::non-urban air or from high stacks",2690.0),
  @("Methane, difluoro-, HFC-32","air::low population density, long-term",2690.0),
  @("Methane, difluoro-, HFC-32","air::lower stratosphere + upper troposphere",2690.0),
  @("Methane, difluoro-, HFC-32","air::unspecified",2690.0),
  @("Methane, fossil","air::urban air close to ground",82.5),
  @("Methane, fossil","air::non-urban air or from high stacks",82.5),
  @("Methane, fossil","air::low population density, long-term",82.5),
  @("Methane, fossil","air::lower stratosphere + upper troposphere",82.5),
  @("Methane, fossil","air::unspecified",82.5),
  @("Methane, from soil or biomass stock","air::non-urban air or from high stacks",82.5),
  @("Methane, from soil or biomass stock","air::low population density, long-term",82.5),
  @("Methane, from soil or biomass stock","air::lower stratosphere + upper troposphere",82.5),
  @("Methane, from soil or biomass stock","air::unspecified",82.5),
  @("Methane, from soil or biomass stock","air::urban air close to ground",82.5),
  @("Methane, monochloro-, R-40","air::urban air close to ground",19.9),
  @("Methane, monochloro-, R-40","air::non-urban air or from high stacks",19.9),
  @("Methane, monochloro-, R-40","air::low population density, long-term",19.9),
  @("Methane, monochloro-, R-40","air::lower stratosphere + upper troposphere",19.9),
  @("Methane, monochloro-, R-40","air::unspecified",19.9),
  @("Methane, non-fossil","air::urban air close to ground",82.5),
  @("Methane, non-fossil","air::non-urban air or from high stacks",82.5),
  @("Methane, non-fossil","air::low population density, long-term",82.5),
  @("Methane, non-fossil","air::lower stratosphere + upper troposphere",82.5),
  @("Methane, non-fossil","air::unspecified",82.5),
  @("Methane, tetrachloro-, R-10","air::urban air close to ground",3810.0),
  @("Methane, tetrachloro-, R-10","air::non-urban air or from high stacks",3810.0),
  @("Methane, tetrachloro-, R-10","air::low population density, long-term",3810.0),
  @("Methane, tetrachloro-, R-10","air::lower stratosphere + upper troposphere",3810.0),
  @("Methane, tetrachloro-, R-10","air::unspecified",3810.0),
  @("Methane, tetrafluoro-, R-14","air::urban air close to ground",5300.0),
  @("Methane, tetrafluoro-, R-14","air::non-urban air or from high stacks",5300.0),
  @("Methane, tetrafluoro-, R-14","air::low population density, long-term",5300.0),
  @("Methane, tetrafluoro-, R-14","air::lower stratosphere + upper troposphere",5300.0),
  @("Methane, tetrafluoro-, R-14","air::unspecified",5300.0),
  @("Methane, trichlorofluoro-, CFC-11","air::urban air close to ground",8320.0),
  @("Methane, trichlorofluoro-, CFC-11","air::non-urban air or from high stacks",8320.0),
  @("Methane, trichlorofluoro-, CFC-11","air::low population density, long-term",8320.0),
  @("Methane, trichlorofluoro-, CFC-11","air::lower stratosphere + upper troposphere",8320.0),
  @("Methane, trichlorofluoro-, CFC-11","air::unspecified",8320.0),
  @("Methane, trifluoro-, HFC-23","air::urban air close to ground",12400.0),
  @("Methane, trifluoro-, HFC-23","air::non-urban air or from high stacks",12400.0),
  @("Methane, trifluoro-, HFC-23","air::low population density, long-term",12400.0),
  @("Methane, trifluoro-, HFC-23","air::lower stratosphere + upper troposphere",12400.0),
  @("Methane, trifluoro-, HFC-23","air::unspecified",12400.0),
  @("Monochloroethane","air::urban air close to ground",1.73),
  @("Monochloroethane","air::non-urban air or from high stacks",1.73),
  @("Monochloroethane","air::low population density, long-term",1.73),
  @("Monochloroethane","air::lower stratosphere + upper troposphere",1.73),
  @("Monochloroethane","air::unspecified",1.73),
  @("Nitric oxide","air::unspecified",16.7),
  @("Nitrogen fluoride","air::urban air close to ground",13400.0),
  @("Nitrogen fluoride","air::non-urban air or from high stacks",13400.0),
  @("Nitrogen fluoride","air::low population density, long-term",13400.0),
  @("Nitrogen fluoride","air::lower stratosphere + upper troposphere",13400.0),
  @("Nitrogen fluoride","air::unspecified",13400.0),
  @("Perfluoropentane","air::unspecified",6680.0),
  @("Propane","air::urban air close to ground",0.072),
  @("Propane","air::non-urban air or from high stacks",0.072),
  @("Propane","air::low population density, long-term",0.072),
  @("Propane","air::lower stratosphere + upper troposphere",0.072),
  @("Propane","air::unspecified",0.072),
  @("Sulfur hexafluoride","air::urban air close to ground",18300.0),
  @("Sulfur hexafluoride","air::non-urban air or from high stacks",18300.0),
  @("Sulfur hexafluoride","air::low population density, long-term",18300.0),
  @("Sulfur hexafluoride","air::lower stratosphere + upper troposphere",18300.0),
  @("Sulfur hexafluoride","air::unspecified",18300.0),
  @("Tetrachloroethylene","air::urban air close to ground",22.8),
  @("Tetrachloroethylene","air::non-urban air or from high stacks",22.8),
  @("Tetrachloroethylene","air::low population density, long-term",22.8),
  @("Tetrachloroethylene","air::lower stratosphere + upper troposphere",22.8),
  @("Tetrachloroethylene","air::unspecified",22.8),
  @("Trichloroethylene","air::urban air close to ground",0.158),
  @("Trichloroethylene","air::non-urban air or from high stacks",0.158),
  @("Trichloroethylene","air::low population density, long-term",0.158),
  @("Trichloroethylene","air::lower stratosphere + upper troposphere",0.158),
  @("Trichloroethylene","air::unspecified",0.158),
  @("VOC, volatile organic compounds","air::unspecified",14.0)
)

$rowCount = $data.Count
$colCount = 3

# Build a 2D array for a single bulk write (row-major -> Excel VARIANT array is [row,col]).
$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i, 0] = $data[$i][0]
    $values[$i, 1] = $data[$i][1]
    $values[$i, 2] = $data[$i][2]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1

# Write the corrected block over the old one.
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 3)).Value = $values

# The corrected table is 4 rows shorter than before (260 -> 256 total rows incl.
# header) since several duplicated/mis-shifted rows collapse away. Remove the
# now-stale trailing rows.
$oldLastRow = 260
if ($oldLastRow -gt $endRow) {
    $ws.Range($ws.Cells.Item($endRow + 1, 1), $ws.Cells.Item($oldLastRow, 3)).EntireRow.Delete() | Out-Null
}

# Refresh the autofilter / filter-database range to match the new extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:C$endRow").AutoFilter() | Out-Null

# Restore the view roughly as last saved.
$ws.Range("A183").Select()
$ws.Range("B196").Select()
